# Fix suspendu mail content issue: add missing "569/SUP 9999" Supervision row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "569/SUP 9999"
$ws.Range("B11").Value = "Supervision"
$ws.Range("C11").Value = "KL365695"
$ws.Range("D11").Value = "YAYA TATA "
$ws.Range("E11").Value = "non"
$ws.Range("F11").Value = "mensuelle"
$ws.Range("G11").Value = 80000
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 15
$ws.Range("J11").Value = 12000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 68000
